# Fix the traceability matrix for R5 ("R5. Assign a vehicle of interest to a client"):
#  - the R5 label was on the wrong row (merged with the tail of R4's last method);
#    move it down one row so it lines up with "Main", like every other requirement block.
#  - add the missing Business-tier methods (printVehicles, three lookForVehicles overloads)
#    and the Vehicle.toString() method that were left out of the matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The "R5..." label in A39 actually belongs on the next row, next to "Main".
$ws.Range("A39").Value = ""
$ws.Range("A40").Value = "R5. Assign a vehicle of interest to a client"

# 2) printVehicles() was missing between addVehicleOfInterest() and the Business tier.
$ws.Rows("41").Insert()
$ws.Range("C41").Value = "printVehicles() : int"

# 3) The three lookForVehicles(...) overloads were missing after searchClient(id : int).
$ws.Rows("45:47").Insert()
$ws.Range("C45").Value = "lookForVehicles(brand : String) : ArrayList<model.Vehicle>"
$ws.Range("C46").Value = "lookForVehicles(model : int) : ArrayList<model.Vehicle>"
$ws.Range("C47").Value = "lookForVehicles(displacement : double) : ArrayList<model.Vehicle>"

# 4) Vehicle.toString() was missing entirely at the end of the R5 block.
$ws.Range("B50").Value = "Vehicle"
$ws.Range("C50").Value = "toString():String"

# Reflect the final editing position/view, like Excel would after scrolling to and
# selecting the last cell typed.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B50").Select()
$excel.ActiveWindow.Zoom = 85
